$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 3: location (client) - update project folder name
$ws.Range("B3").Value = "W:\Projects\תכניות מרחביות\בדיקה מרחבית גילה\קבצי עבודה\תחזיות_דמוגרפיות"

# Row 4: forecast_version - "max" -> "realy"
$ws.Range("B4").Value = "realy"

# Row 5: v_date - "240410_with_poten" -> "240818_with_poten"
$ws.Range("B5").Value = "240818_with_poten"
